$wb = $excel.ActiveWorkbook

# Rename the "nome" sheet to "nomes"
$ws = $wb.Worksheets.Item("nome")
$ws.Name = "nomes"

# Fix header text and add new column/row of data
$ws.Range("A1").Value = "Nome"
$ws.Range("B1").Value = "Idade"
$ws.Range("C1").Value = "Sexo"

$ws.Range("A2").Value = "carlos"
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "41"
$ws.Range("C2").Value = "M"
